$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New Price (column D) and Volume(1h) (column E) values for rows 2-51,
# taken from the updated cryptos list.
$data = @(
    @{ Row = 2; D = "27.449.45"; E = "  -2.24%  " },
    @{ Row = 3; D = "1.836.97"; E = "  -2.71%  " },
    @{ Row = 4; D = "1.003"; E = "  -1.03%  " },
    @{ Row = 5; D = "332.98"; E = "  -0.94%  " },
    @{ Row = 6; D = "1.003"; E = "  -0.90%  " },
    @{ Row = 7; D = "0.4615"; E = "  -3.09%  " },
    @{ Row = 8; D = "0.3805"; E = "  -3.79%  " },
    @{ Row = 9; D = "46.60"; E = "  -1.11%  " },
    @{ Row = 10; D = $null; E = "  -2.06%  " },
    @{ Row = 11; D = "0.9726"; E = "  -4.70%  " },
    @{ Row = 12; D = "21.06"; E = "  -3.91%  " },
    @{ Row = 13; D = "1.861.48"; E = "  -2.42%  " },
    @{ Row = 14; D = "5.888"; E = "  -2.67%  " },
    @{ Row = 15; D = "7.026"; E = "  -2.63%  " },
    @{ Row = 16; D = "1.004"; E = "  -1.07%  " },
    @{ Row = 17; D = "87.83"; E = "  -0.87%  " },
    @{ Row = 18; D = "0.06640"; E = "  -1.78%  " },
    @{ Row = 19; D = $null; E = "  -2.26%  " },
    @{ Row = 20; D = "16.98"; E = "  -0.66%  " },
    @{ Row = 21; D = "1.004"; E = "  -0.78%  " },
    @{ Row = 22; D = "27.444.01"; E = "  -2.17%  " },
    @{ Row = 23; D = "5.336"; E = "  -3.36%  " },
    @{ Row = 24; D = "10.83"; E = "  -1.88%  " },
    @{ Row = 25; D = "2.302"; E = "  -2.16%  " },
    @{ Row = 26; D = "157.66"; E = "  -1.29%  " },
    @{ Row = 27; D = "19.35"; E = "  -3.35%  " },
    @{ Row = 28; D = "2.065"; E = "  -2.06%  " },
    @{ Row = 29; D = "5.316"; E = "  -3.74%  " },
    @{ Row = 30; D = "118.80"; E = "  -2.52%  " },
    @{ Row = 31; D = "0.9515"; E = "  -2.75%  " },
    @{ Row = 32; D = "0.09300"; E = "  -3.09%  " },
    @{ Row = 33; D = "3.566"; E = "  -1.99%  " },
    @{ Row = 34; D = "5.226"; E = "  -2.58%  " },
    @{ Row = 35; D = "1.320"; E = "  -3.39%  " },
    @{ Row = 36; D = "0.05931"; E = "  -2.39%  " },
    @{ Row = 37; D = "0.02191"; E = "  -2.85%  " },
    @{ Row = 38; D = "8.070"; E = "  -1.81%  " },
    @{ Row = 39; D = "1.160"; E = "  -4.06%  " },
    @{ Row = 40; D = "0.5784"; E = "  -3.35%  " },
    @{ Row = 41; D = "0.1837"; E = "  -3.05%  " },
    @{ Row = 42; D = "10.01"; E = "  -3.38%  " },
    @{ Row = 43; D = "1.241"; E = "  -1.75%  " },
    @{ Row = 44; D = "0.5483"; E = "  -3.27%  " },
    @{ Row = 45; D = "11.92"; E = "  -2.65%  " },
    @{ Row = 46; D = "1.864"; E = "  -3.70%  " },
    @{ Row = 47; D = "0.06653"; E = "  -2.54%  " },
    @{ Row = 48; D = "109.91"; E = "  -2.20%  " },
    @{ Row = 49; D = "1.041"; E = "  -3.01%  " },
    @{ Row = 50; D = "1.003"; E = "  -0.99%  " },
    @{ Row = 51; D = "69.51"; E = "  -2.23%  " }
)

foreach ($item in $data) {
    $r = $item.Row
    if ($item.D -ne $null) {
        $dCell = $ws.Range("D$r")
        # Force a Text format before assigning so decimal-looking strings
        # (e.g. "10.01") are stored as text instead of being converted to numbers.
        $dCell.NumberFormat = "@"
        $dCell.Value = $item.D
    }
    $ws.Range("E$r").Value = $item.E
}

# Restore the original (default) cell style on the Price cells that were
# temporarily switched to Text format, so no formatting change is introduced.
foreach ($item in $data) {
    if ($item.D -ne $null) {
        $ws.Range("D" + $item.Row).Style = "Normal"
    }
}

